$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure Price column (D) is treated as text so values like "1.003" or "241.78"
# are not reinterpreted as numbers by Excel, matching the original inline-string cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.499.90'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '1.858.13'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').Value = '241.78'
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('D6').Value = '0.6343'
$ws.Range('E6').Value = '  +1.28%  '
$ws.Range('D7').Value = '1.004'
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('D8').Value = '0.07582'
$ws.Range('E8').Value = '  -1.22%  '
$ws.Range('D9').Value = '0.2929'
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('D10').Value = '24.62'
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('D11').Value = '0.07771'
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').Value = '1.858.28'
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('D13').Value = '5.046'
$ws.Range('E13').Value = '  +0.39%  '
$ws.Range('D14').Value = '0.6865'
$ws.Range('E14').Value = '  +0.78%  '
$ws.Range('E15').Value = '  -1.95%  '
$ws.Range('D16').Value = '83.51'
$ws.Range('D17').Value = '2.119.15'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D19').Value = '29.521.28'
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').Value = '230.37'
$ws.Range('E20').Value = '  +0.98%  '
$ws.Range('E21').Value = '  +0.17%  '
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('D23').Value = '7.515'
$ws.Range('E23').Value = '  +1.42%  '
$ws.Range('D24').Value = '1.005'
$ws.Range('E24').Value = '  +0.34%  '
$ws.Range('D25').Value = '159.60'
$ws.Range('E25').Value = '  +1.34%  '
$ws.Range('D26').Value = '0.1395'
$ws.Range('E26').Value = '  +1.64%  '
$ws.Range('D27').Value = '8.479'
$ws.Range('E27').Value = '  +1.02%  '
$ws.Range('D28').Value = '17.75'
$ws.Range('E28').Value = '  +0.38%  '
$ws.Range('D29').Value = '1.424'
$ws.Range('E29').Value = '  +6.09%  '
$ws.Range('D30').Value = '1.480'
$ws.Range('E30').Value = '  +0.96%  '
$ws.Range('D31').Value = '0.05705'
$ws.Range('E31').Value = '  +1.16%  '
$ws.Range('E32').Value = '  +1.01%  '
$ws.Range('D33').Value = '4.073'
$ws.Range('E33').Value = '  +1.02%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '1.160'
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = '1.830'
$ws.Range('E35').Value = '  -0.57%  '
$ws.Range('D36').Value = '0.6977'
$ws.Range('E36').Value = '  -1.47%  '
$ws.Range('D37').Value = '2.596'
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = '1.257.64'
$ws.Range('E38').Value = '  +2.28%  '
$ws.Range('D39').Value = '0.01838'
$ws.Range('E39').Value = '  +3.00%  '
$ws.Range('D40').Value = '2.787'
$ws.Range('E40').Value = '  +0.68%  '
$ws.Range('D41').Value = '6.512'
$ws.Range('E41').Value = '  -0.48%  '
$ws.Range('D42').Value = '0.9092'
$ws.Range('E42').Value = '  +0.62%  '
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('D44').Value = '2.020.96'
$ws.Range('E44').Value = '  -0.62%  '
$ws.Range('D45').Value = '101.66'
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('D46').Value = '66.30'
$ws.Range('E46').Value = '  +0.71%  '
$ws.Range('D47').Value = '7.147'
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('D48').Value = '0.1171'
$ws.Range('E48').Value = '  +1.83%  '
$ws.Range('D49').Value = '9.051'
$ws.Range('E49').Value = '  +0.92%  '
$ws.Range('D50').Value = '0.3984'
$ws.Range('E50').Value = '  -0.71%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.00000000115'
$ws.Range('E51').Value = '  -4.66%  '
